# Update "想去人数" (interest count) figures on the gh-pages regenerated
# output, as of commit 456a3b4.
#
# Sheet "展览" (Exhibitions) rows 3,5,12,13,14,15 and the aggregate sheet
# "全部类型" (All types) rows 4,6,14,15,16,17 hold the same six events;
# column F is the count that changed.

$wb = $excel.ActiveWorkbook

$updates = @{
    "展览"   = @{ 3 = 1073; 5 = 2927; 12 = 156; 13 = 72; 14 = 2764; 15 = 1039 }
    "全部类型" = @{ 4 = 1073; 6 = 2927; 14 = 156; 15 = 72; 16 = 2764; 17 = 1039 }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Cells.Item($row, 6).Value = $rows[$row]
    }
}
